$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-5 (Sending cluster changed for each row because the
# underlying per-cluster TPM values were recomputed; Target cluster labels
# stay the same but all derived statistics are refreshed with new numbers).

$rows = @(
    @{ Row = 2;  A = "FAPs";  D = "ECs";
       E = 2;  F = 0.6666666666666666;      G = 0.344207;              H = 1.032621;
       I = 0.6985282229833164;  J = 0.6985282229833165;
       K = 1;  L = 0.3333333333333333;
       M = 0.006923666666666667; N = 0.020771;
       O = 0.01563438526027703;  P = 0.01563438526027703;
       Q = 0.002383174532333333; R = 0.021448570791;
       S = 0.01092105935329787;  T = 0.01092105935329787 },

    @{ Row = 3;  A = "FAPs";  D = "MuSCs";
       E = 2;  F = 0.6666666666666666;      G = 0.344207;              H = 1.032621;
       I = 0.6985282229833164;  J = 0.6985282229833165;
       K = 3;  L = 1;
       M = 0.435925;              N = 1.307775;
       O = 0.9843656147397229;   P = 0.9843656147397229;
       Q = 0.150048436475;       R = 1.350435928275;
       S = 0.6876071636300185;   T = 0.6876071636300186 },

    @{ Row = 4;  A = "MuSCs"; D = "ECs";
       E = 1;  F = 0.3333333333333333;      G = 0.1485533333333333;    H = 0.44566;
       I = 0.3014717770166836;  J = 0.3014717770166836;
       K = 1;  L = 0.3333333333333333;
       M = 0.006923666666666667; N = 0.020771;
       O = 0.01563438526027703;  P = 0.01563438526027703;
       Q = 0.001028533762222222; R = 0.009256803860000001;
       S = 0.004713325906979161; T = 0.00471332590697916 },

    @{ Row = 5;  A = "MuSCs"; D = "MuSCs";
       E = 1;  F = 0.3333333333333333;      G = 0.1485533333333333;    H = 0.44566;
       I = 0.3014717770166836;  J = 0.3014717770166836;
       K = 3;  L = 1;
       M = 0.435925;              N = 1.307775;
       O = 0.9843656147397229;   P = 0.9843656147397229;
       Q = 0.06475811183333333;  R = 0.5828230065;
       S = 0.2967584511097044;   T = 0.2967584511097044 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A   # A: Sending cluster
    $ws.Cells.Item($row, 4).Value = $r.D   # D: Target cluster

    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
